# Applies the "cleaning up the data so its useable" edit:
#   1. Bump the auto date field shown on every slide layout + the slide
#      master from 1/12/19 to 1/13/19.
#   2. Merge the "Total " + "Number of " runs in the "Total Number of
#      schools by size" callout into a single "Total Number of " run.
#   3. Replace the "$40,000" callout text with "250".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder: 1/12/19 -> 1/13/19 on the slide master and on
#    every custom (slide) layout.
# ---------------------------------------------------------------------
function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq "1/12/19") {
                    $shp.TextFrame.TextRange.Text = "1/13/19"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    Update-DatePlaceholder $master.CustomLayouts.Item($li)
}

# ---------------------------------------------------------------------
# 2) "Total " / "Number of " runs -> single "Total Number of " run.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 4" -and $shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -like "Total Number of schools*") {
            $tr.Characters(1, 16).Text = "Total Number of "
        }
    }
}

# ---------------------------------------------------------------------
# 3) "$40,000" callout -> "250".
# ---------------------------------------------------------------------
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Name -eq "Rectangular Callout 100" -and $shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "`$40,000") {
            $tr.Characters(1, 1).Text = "2"
            $tr.Characters(2, 1).Text = "5"
            $tr.Characters(3, 1).Text = "0"
            $tr.Characters(4, 4).Text = ""
        }
    }
}
